$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.03273466666666667
$ws.Range("H2").Value = 0.098204
$ws.Range("I2").Value = 0.08359843399780884
$ws.Range("J2").Value = 0.08359843399780884
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 0.01403475919066667
$ws.Range("R2").Value = 0.126312832716
$ws.Range("S2").Value = 0.0003444676773873962
$ws.Range("T2").Value = 0.0003444676773873962
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.03273466666666667
$ws.Range("H3").Value = 0.098204
$ws.Range("I3").Value = 0.08359843399780884
$ws.Range("J3").Value = 0.08359843399780884
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 2.626179126115999
$ws.Range("R3").Value = 23.635612135044
$ws.Range("S3").Value = 0.06445666873843021
$ws.Range("T3").Value = 0.06445666873843021
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.03273466666666667
$ws.Range("H4").Value = 0.098204
$ws.Range("I4").Value = 0.08359843399780884
$ws.Range("J4").Value = 0.08359843399780884
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 0.7658644404591112
$ws.Range("R4").Value = 6.892779964132
$ws.Range("S4").Value = 0.01879729758199123
$ws.Range("T4").Value = 0.01879729758199123
$ws.Range("I5").Value = 0.3399848984133119
$ws.Range("J5").Value = 0.3399848984133119
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 0.05707769810400001
$ws.Range("R5").Value = 0.5136992829360001
$ws.Range("S5").Value = 0.001400909116387192
$ws.Range("T5").Value = 0.001400909116387193
$ws.Range("I6").Value = 0.3399848984133119
$ws.Range("J6").Value = 0.3399848984133119
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("R6").Value = 96.123226314024
$ws.Range("S6").Value = 0.2621376134111564
$ws.Range("T6").Value = 0.2621376134111565
$ws.Range("I7").Value = 0.3399848984133119
$ws.Range("J7").Value = 0.3399848984133119
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 3.114679684008
$ws.Range("R7").Value = 28.032117156072
$ws.Range("S7").Value = 0.07644637588576826
$ws.Range("T7").Value = 0.07644637588576826
$ws.Range("E8").Value = 2.0
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.2257076666666666
$ws.Range("H8").Value = 0.6771229999999999
$ws.Range("I8").Value = 0.5764166675888793
$ws.Range("J8").Value = 0.5764166675888793
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 0.09677058212966666
$ws.Range("R8").Value = 0.8709352391669999
$ws.Range("S8").Value = 0.002375127154857092
$ws.Range("T8").Value = 0.002375127154857092
$ws.Range("E9").Value = 2.0
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.2257076666666666
$ws.Range("H9").Value = 0.6771229999999999
$ws.Range("I9").Value = 0.5764166675888793
$ws.Range("J9").Value = 0.5764166675888793
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("N9").Value = 240.678711
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 18.10767675871699
$ws.Range("R9").Value = 162.9690908284529
$ws.Range("S9").Value = 0.4444329447494203
$ws.Range("T9").Value = 0.4444329447494203
$ws.Range("E10").Value = 2.0
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.2257076666666666
$ws.Range("H10").Value = 0.6771229999999999
$ws.Range("I10").Value = 0.5764166675888793
$ws.Range("J10").Value = 0.5764166675888793
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 5.280685384678778
$ws.Range("R10").Value = 47.526168462109
$ws.Range("S10").Value = 0.129608595684602
$ws.Range("T10").Value = 0.129608595684602
